$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 20279.8
$ws.Range("I12").Value = 20279.8
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 20279.8
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -20109.8

$ws.Range("H33").Value = 859
$ws.Range("I33").Value = 406.33334
$ws.Range("J33").Value = 1476.2727
$ws.Range("K33").Value = 406.33334
$ws.Range("L33").Value = 1476.2727
$ws.Range("M33").Value = -177.33334
$ws.Range("N33").Value = -1934.2727

$ws.Range("H41").Value = 1066.9584
$ws.Range("I41").Value = 1239.3
$ws.Range("J41").Value = 943.8570999999999
$ws.Range("K41").Value = 1239.3
$ws.Range("L41").Value = 943.8570999999999
$ws.Range("M41").Value = -799.3
$ws.Range("N41").Value = -1823.8571

$ws.Range("H70").Value = 11892.786
$ws.Range("I70").Value = 2249.8333
$ws.Range("J70").Value = 19125
$ws.Range("K70").Value = 6749.499899999999
$ws.Range("L70").Value = 57375
$ws.Range("M70").Value = -6479.499899999999

$ws.Range("H73").Value = 11892.786
$ws.Range("I73").Value = 2249.8333
$ws.Range("J73").Value = 19125
$ws.Range("K73").Value = 6749.499899999999
$ws.Range("L73").Value = 57375
$ws.Range("M73").Value = -5813.499899999999

$ws.Range("H101").Value = 533.8461
$ws.Range("I101").Value = 937.75
$ws.Range("J101").Value = 354.33334
$ws.Range("K101").Value = 2813.25
$ws.Range("L101").Value = 1063.00002
$ws.Range("M101").Value = -1191.25
$ws.Range("N101").Value = -4307.000019999999

$ws.Range("H121").Value = 1866.6666
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1866.6666
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5599.9998
$ws.Range("N121").Value = -9093.9998

$ws.Range("H132").Value = 3581.4285
$ws.Range("I132").Value = 3820.077
$ws.Range("J132").Value = 479
$ws.Range("K132").Value = 11460.231
$ws.Range("L132").Value = 1437
$ws.Range("M132").Value = -8930.231
$ws.Range("N132").Value = -6497

$ws.Range("H135").Value = 23810128
$ws.Range("I135").Value = 23810128
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 214291152
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -214288617

$ws.Range("H137").Value = 2296.3333
$ws.Range("I137").Value = 1967.44
$ws.Range("J137").Value = 3043.818
$ws.Range("K137").Value = 5902.32
$ws.Range("L137").Value = 9131.454000000002
$ws.Range("M137").Value = -3352.32
$ws.Range("N137").Value = -14231.454

$ws.Range("H138").Value = 2244.52
$ws.Range("I138").Value = 1177.6571
$ws.Range("J138").Value = 3178.025
$ws.Range("K138").Value = 3532.9713
$ws.Range("L138").Value = 9534.075000000001
$ws.Range("M138").Value = 1607.0287
$ws.Range("N138").Value = -19814.075

$ws.Range("H141").Value = 1001.7692
$ws.Range("I141").Value = 1001.7692
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3005.3076
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2174.6924
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2285.1396
$ws.Range("I32").Value = 2129.394
$ws.Range("J32").Value = 2799.1
$ws.Range("K32").Value = 2129.394
$ws.Range("L32").Value = 2799.1
$ws.Range("M32").Value = -1842.394
$ws.Range("N32").Value = -3373.1

$ws.Range("H74").Value = 27030908
$ws.Range("I74").Value = 30306914
$ws.Range("J74").Value = 3847.25
$ws.Range("K74").Value = 30306914
$ws.Range("L74").Value = 3847.25
$ws.Range("M74").Value = -30306040
$ws.Range("N74").Value = -5595.25

$ws.Range("H77").Value = 27030908
$ws.Range("I77").Value = 30306914
$ws.Range("J77").Value = 3847.25
$ws.Range("K77").Value = 151534570
$ws.Range("L77").Value = 19236.25
$ws.Range("M77").Value = -151530202
$ws.Range("N77").Value = -27972.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2262.4285
$ws.Range("I105").Value = 2262.4285
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2262.4285
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -515.4285
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 26843144
$ws.Range("I134").Value = 34000930
$ws.Range("J134").Value = 1449.25
$ws.Range("K134").Value = 102002790
$ws.Range("L134").Value = 4347.75
$ws.Range("M134").Value = -102000255
$ws.Range("N134").Value = -9417.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9189.581
$ws.Range("I31").Value = 6497.9375
$ws.Range("J31").Value = 12060.667
$ws.Range("K31").Value = 6497.9375
$ws.Range("L31").Value = 12060.667
$ws.Range("M31").Value = -6202.9375
$ws.Range("N31").Value = -12650.667

$ws.Range("H34").Value = 9189.581
$ws.Range("I34").Value = 6497.9375
$ws.Range("J34").Value = 12060.667
$ws.Range("K34").Value = 6497.9375
$ws.Range("L34").Value = 12060.667
$ws.Range("M34").Value = -6295.9375
$ws.Range("N34").Value = -12464.667

$ws.Range("H58").Value = 15155736
$ws.Range("I58").Value = 20004976
$ws.Range("J58").Value = 1863.625
$ws.Range("K58").Value = 20004976
$ws.Range("L58").Value = 1863.625
$ws.Range("M58").Value = -20004773
$ws.Range("N58").Value = -2269.625

$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240
$ws.Range("M65").ClearContents()

$ws.Range("H136").Value = 15155736
$ws.Range("I136").Value = 20004976
$ws.Range("J136").Value = 1863.625
$ws.Range("K136").Value = 60014928
$ws.Range("L136").Value = 5590.875
$ws.Range("M136").Value = -60012378
$ws.Range("N136").Value = -10690.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 100681.5
$ws.Range("I5").Value = 143263.58
$ws.Range("J5").Value = 1323.3334
$ws.Range("K5").Value = 429790.74
$ws.Range("L5").Value = 3970.0002
$ws.Range("M5").Value = -429678.74

$ws.Range("H56").Value = 15560.209
$ws.Range("I56").Value = 15560.209
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 15560.209
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -15030.209

$ws.Range("H68").Value = 4999.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4999.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14998.5
$ws.Range("N68").Value = -16620.5

$ws.Range("H71").Value = 4999.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4999.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 44995.5
$ws.Range("N71").Value = -53107.5

$ws.Range("H103").Value = 997.2
$ws.Range("I103").Value = 194.83333
$ws.Range("J103").Value = 1532.1111
$ws.Range("K103").Value = 584.49999
$ws.Range("L103").Value = 4596.3333
$ws.Range("M103").Value = 294.50001
$ws.Range("N103").Value = -6354.3333

$ws.Range("H121").Value = 881936
$ws.Range("I121").Value = 113612.555
$ws.Range("J121").Value = 2264918.2
$ws.Range("K121").Value = 340837.665
$ws.Range("L121").Value = 6794754.600000001
$ws.Range("M121").Value = -339527.665
$ws.Range("N121").Value = -6797374.600000001

$ws.Range("H132").Value = 1399.2858
$ws.Range("I132").Value = 1215.8334
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10942.5006
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -8412.500599999999
$ws.Range("N132").Value = -27560

$ws.Range("H135").Value = 100681.5
$ws.Range("I135").Value = 143263.58
$ws.Range("J135").Value = 1323.3334
$ws.Range("K135").Value = 1289372.22
$ws.Range("L135").Value = 11910.0006
$ws.Range("M135").Value = -1286837.22

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4728.486
$ws.Range("I122").Value = 3109.2
$ws.Range("J122").Value = 8776.700000000001
$ws.Range("K122").Value = 9327.599999999999
$ws.Range("L122").Value = 26330.1
$ws.Range("M122").Value = -6877.599999999999

$ws.Range("H132").Value = 13891927
$ws.Range("I132").Value = 13891927
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 41675781
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -41673251
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2821.3447
$ws.Range("I16").Value = 1479.8334
$ws.Range("J16").Value = 3768.2942
$ws.Range("K16").Value = 1479.8334
$ws.Range("L16").Value = 3768.2942
$ws.Range("M16").Value = -1309.8334
$ws.Range("N16").Value = -4108.2942

$ws.Range("H22").Value = 5766.6665
$ws.Range("I22").Value = 5025
$ws.Range("J22").Value = 7250
$ws.Range("K22").Value = 5025
$ws.Range("L22").Value = 7250
$ws.Range("M22").Value = -4730

$ws.Range("H27").Value = 5766.6665
$ws.Range("I27").Value = 5025
$ws.Range("J27").Value = 7250
$ws.Range("K27").Value = 5025
$ws.Range("L27").Value = 7250
$ws.Range("M27").Value = -4918

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 4999.5
$ws.Range("I122").Value = 4999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12548.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1228.1177
$ws.Range("I113").Value = 1282.3334
$ws.Range("J113").Value = 821.5
$ws.Range("K113").Value = 3847.0002
$ws.Range("L113").Value = 2464.5
$ws.Range("M113").Value = -1677.0002
$ws.Range("N113").Value = -6804.5
